# Fix Training Data Issue (#48)
# The BF column holds a "Date" label/value that was mistakenly written as
# "5-11-2007-08" (a mangled mashup of the file name / season). It should
# read the actual game date "2008-05-11" (ISO yyyy-mm-dd) for every data
# row (rows 2-31; row 1 is the "Date" header).
#
# NumberFormat is forced to Text ("@") before the assignment so Excel does
# not reinterpret the date-shaped string as a real date serial, then the
# cell Style is reset back to "Normal" so no stray formatting/style is left
# behind on the cell (keeps the cell's style identical to before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-11-2007-08"
$newValue = "2008-05-11"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
